# Remove Wong, add Estapa to flux comparison figure
#
# 1. Rename the existing "to_df" sheet to "thorium".
# 2. Add a new "traps" sheet (right after "thorium") holding the Estapa
#    sediment-trap flux data (depth, flux, flux_u), and make it the active
#    sheet/tab.

$wb = $excel.ActiveWorkbook

# --- rename existing sheet -------------------------------------------------
$thorium = $wb.Worksheets.Item(1)
$thorium.Name = "thorium"

# --- add the new "traps" sheet right after "thorium" ------------------------
$traps = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $thorium)
$traps.Name = "traps"

# --- headers -----------------------------------------------------------
$traps.Range("A1").Value = "depth"
$traps.Range("B1").Value = "flux"
$traps.Range("C1").Value = "flux_u"

# --- data rows -----------------------------------------------------------
$data = @(
    @(100, 1.3848625432886685,  0.7719548047810465),
    @(150, 1.0407822367945465,  0.38784151665860167),
    @(200, 0.75568080420966999, 0.23729528471901881),
    @(330, 0.64319121464272433, 0.28415347289088116),
    @(500, 0.85783057996062873, 0.43077789507666842)
)

$r = 2
foreach ($row in $data) {
    $traps.Cells.Item($r, 1).Value = $row[0]
    $traps.Cells.Item($r, 2).Value = $row[1]
    $traps.Cells.Item($r, 3).Value = $row[2]
    $r++
}

# --- selection / active sheet -----------------------------------------------
$traps.Range("E13").Select() | Out-Null
$traps.Activate() | Out-Null
